$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply updated odds values to row 3 and row 4 as per the 2024-10-10 FlashScore data refresh
$ws.Range("G3").Value = 1.53
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 6.5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("W3").Value = 5
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 101
$ws.Range("AK3").Value = 81
$ws.Range("AM3").Value = 67
$ws.Range("AW3").Value = 8
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 2.18
$ws.Range("J4").Value = 3.45
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 2.72
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 3.4
$ws.Range("R4").Value = 1.98
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.14
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.12
$ws.Range("W4").Value = 11
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 37
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 6.8
$ws.Range("AE4").Value = 12.5
$ws.Range("AH4").Value = 9.25
$ws.Range("AI4").Value = 11.75
$ws.Range("AJ4").Value = 8.75
$ws.Range("AK4").Value = 22
$ws.Range("AL4").Value = 16
$ws.Range("AM4").Value = 23
$ws.Range("AN4").Value = 5
$ws.Range("AO4").Value = 15.5
$ws.Range("AT4").Value = 2.92
$ws.Range("AU4").Value = 6.6
$ws.Range("AW4").Value = 4.2
$ws.Range("AX4").Value = 11
$ws.Range("AY4").Value = 17
$ws.Range("AZ4").Value = 40
$ws.Range("BA4").Value = 65
$ws.Range("BB4").Value = 200
